# Adds a new "Address" column (F) to the sheet, shifting the existing
# "District" column (and everything to its right) one column to the right.
# The new Address values are derived from the school/location portion of
# each teacher's entry in column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at F; old F (District) and beyond shift to G, etc.
$ws.Columns.Item(6).Insert()

# New header for the inserted column.
$ws.Range("F2").Value = "Address"

# Map of row number -> Address text to place in column F.
$addresses = @{
    3 = 'S M H S LakkalakattiRon'
    4 = 'J T Jr Comp. CollegeDambal Mundaragi'
    5 = 'K S S H Gajendragad'
    6 = 'G H S Chikkamandya'
    7 = 'M S Dambal Girls High SchoolMundargi'
    8 = 'S M High School Balaganur'
    9 = 'G H S GojanurShirahatti'
    10 = 'G H S MagadiShirahatti'
    11 = 'Sanjay High School HolealurRon'
    12 = 'S R E S High School SuranagiShirahatti'
    13 = 'S A V H S YarebeleriRon'
    14 = 'S S B H Ron'
    15 = 'G H S ValageremenasaK R Pet'
    16 = 'S H B H S High School MadalagariRon'
    17 = 'Govt. High School Thimmana Hosur'
    19 = 'U G H S (RMSA) DindurRon'
    20 = 'G H S KadiluvagiluMaddur'
    21 = 'G H S Kalakapur (RMSA)Ron'
    22 = 'G H S Soonagahalli'
    23 = 'G H S Papanashi'
    24 = 'H C E S High School Chinchali'
    25 = 'Model High School ShanatagiriRon'
    26 = 'Nargund'
    27 = 'G H S HalehalliMaddur'
    28 = 'G H S Goravale Mandya North'
    29 = 'G H S (RMSA) ChikkamannurRon'
    30 = 'G H S MugaliRon'
    31 = 'G H S Chujjala KyathanahallyK R Pet'
    32 = 'Govt Junior College K R Pet'
    33 = 'Sharadha Girls High SchoolRon'
    34 = 'S A V H S NaregalRon'
    35 = 'S S B H S NagendragadRon'
    36 = 'S S H LaxmeshawarShiarahatti'
    37 = 'Shri Sangameshwar High SchoolSavadiRon'
    38 = 'S A H S NidagundiRon'
    39 = 'Sri Jagadguru AnnadaneshwarHigh School – Mundargi'
    41 = 'M D R School HirewaddattiMundaragi'
    42 = 'Govt. High School Satanuru'
    43 = 'Govt. Junior College Besagarahalli Maddur'
    44 = 'G H P S KonchigeriShirahatti'
    45 = 'G H S HalakereRon'
}

foreach ($row in $addresses.Keys) {
    $ws.Range("F$row").Value = $addresses[$row]
}
